$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = '@'
    $rng.Value = $text
    $rng.ClearFormats()
}

Set-TextValue 'D2' '30.437.42'
Set-TextValue 'E2' '  +0.82%  '
Set-TextValue 'D3' '1.868.90'
Set-TextValue 'E3' '  +0.32%  '
Set-TextValue 'D4' '1.000'
Set-TextValue 'E4' '  +0.11%  '
Set-TextValue 'D5' '246.35'
Set-TextValue 'E5' '  +1.39%  '
Set-TextValue 'D6' '1.000'
Set-TextValue 'E6' '  +0.08%  '
Set-TextValue 'D7' '0.4740'
Set-TextValue 'E7' '  +0.43%  '
Set-TextValue 'D8' '0.2903'
Set-TextValue 'E8' '  +1.57%  '
Set-TextValue 'D9' '0.06496'
Set-TextValue 'E9' '  +0.27%  '
Set-TextValue 'D10' '21.98'
Set-TextValue 'E10' '  +5.66%  '
Set-TextValue 'B11' 'Litecoin'
Set-TextValue 'C11' 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue 'D11' '97.84'
Set-TextValue 'E11' '  +3.92%  '
Set-TextValue 'B12' 'TRON'
Set-TextValue 'C12' 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue 'D12' '0.07711'
Set-TextValue 'E12' '  +0.40%  '
Set-TextValue 'D13' '0.7359'
Set-TextValue 'E13' '  +7.83%  '
Set-TextValue 'D14' '1.869.61'
Set-TextValue 'E14' '  +0.51%  '
Set-TextValue 'D15' '5.117'
Set-TextValue 'E15' '  +0.79%  '
Set-TextValue 'D16' '274.14'
Set-TextValue 'E16' '  +1.61%  '
Set-TextValue 'D17' '30.402.71'
Set-TextValue 'E17' '  +0.74%  '
Set-TextValue 'D18' '13.37'
Set-TextValue 'E18' '  +0.11%  '
Set-TextValue 'D19' '0.000007562'
Set-TextValue 'E19' '  +0.21%  '
Set-TextValue 'D20' '1.001'
Set-TextValue 'E20' '  +0.13%  '
Set-TextValue 'D21' '2.117.17'
Set-TextValue 'E21' '  +0.53%  '
Set-TextValue 'E22' '  +0.19%  '
Set-TextValue 'D23' '5.230'
Set-TextValue 'E23' '  +1.04%  '
Set-TextValue 'D24' '6.166'
Set-TextValue 'E24' '  +1.06%  '
Set-TextValue 'D25' '9.263'
Set-TextValue 'E25' '  -0.88%  '
Set-TextValue 'D26' '164.31'
Set-TextValue 'E26' '  -1.03%  '
Set-TextValue 'D27' '18.83'
Set-TextValue 'E27' '  +0.53%  '
Set-TextValue 'D28' '1.925'
Set-TextValue 'E28' '  +2.15%  '
Set-TextValue 'E29' '  +1.68%  '
Set-TextValue 'D30' '1.367'
Set-TextValue 'E30' '  -0.66%  '
Set-TextValue 'D31' '1.506'
Set-TextValue 'E31' '  -0.06%  '
Set-TextValue 'D32' '4.306'
Set-TextValue 'E32' '  +1.68%  '
Set-TextValue 'D33' '4.162'
Set-TextValue 'E33' '  +4.11%  '
Set-TextValue 'D34' '0.04826'
Set-TextValue 'E34' '  +2.49%  '
Set-TextValue 'D35' '1.118'
Set-TextValue 'E35' '  +0.64%  '
Set-TextValue 'D36' '0.6967'
Set-TextValue 'E36' '  +1.52%  '
Set-TextValue 'E37' '  +0.14%  '
Set-TextValue 'D38' '0.01854'
Set-TextValue 'E38' '  +1.41%  '
Set-TextValue 'D39' '2.745'
Set-TextValue 'E39' '  +0.69%  '
Set-TextValue 'D40' '6.307'
Set-TextValue 'E40' '  -1.25%  '
Set-TextValue 'B41' 'Aave'
Set-TextValue 'C41' 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue 'D41' '72.44'
Set-TextValue 'E41' '  +3.14%  '
Set-TextValue 'B42' 'RenderToken'
Set-TextValue 'C42' 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 'D42' '1.968'
Set-TextValue 'E42' '  +4.10%  '
Set-TextValue 'D43' '0.4185'
Set-TextValue 'E43' '  +2.77%  '
Set-TextValue 'D44' '1.000'
Set-TextValue 'E44' '  +0.11%  '
Set-TextValue 'D45' '0.8347'
Set-TextValue 'E45' '  -0.05%  '
Set-TextValue 'D46' '102.74'
Set-TextValue 'E46' '  +0.59%  '
Set-TextValue 'D47' '9.257'
Set-TextValue 'E47' '  +0.10%  '
Set-TextValue 'D48' '7.020'
Set-TextValue 'E48' '  +1.10%  '
Set-TextValue 'D49' '35.36'
Set-TextValue 'E49' '  +2.73%  '
Set-TextValue 'D50' '920.62'
Set-TextValue 'E50' '  -0.71%  '
Set-TextValue 'D51' '0.05634'
Set-TextValue 'E51' '  +1.42%  '
